# "set property's private value as true"
# Column D on the "Property" sheet holds the boolean "Private" field; the
# author flipped it to TRUE for every data row (rows 2-6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D6").Value = $true

# The edited column becomes the sheet's active selection (mirrors the
# author re-selecting D2:D6 after the edit, in place of the old F2:F6).
[void]$ws.Range("D2:D6").Select()

# Column F already carries a TRUE/FALSE list validation (split across
# F2:F6 and F7:F1048576); normalize it to one contiguous area and extend
# the same TRUE/FALSE dropdown validation to the newly-edited column D.
$fRange = $ws.Range("F2:F1048576")
$fRange.Validation.Delete()
$fRange.Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$fRange.Validation.IgnoreBlank = $true
$fRange.Validation.InCellDropdown = $true
$fRange.Validation.ShowInput = $true
$fRange.Validation.ShowError = $true

$dRange = $ws.Range("D2:D6")
$dRange.Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$dRange.Validation.IgnoreBlank = $true
$dRange.Validation.InCellDropdown = $true
$dRange.Validation.ShowInput = $true
$dRange.Validation.ShowError = $true

$wb.Save()
